$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect so values can be written
$ws.Unprotect()

# Update the confidential footer text (date bump 2021-05-12 -> 2021-05-13)
$oldText = $ws.Range("A59").Value2
$newText = $oldText -replace "2021-05-12", "2021-05-13"
$ws.Range("A59").Value2 = $newText

# Update Weight (D) and Percent Change (E) values for holdings rows 2-56
$ws.Range("D2").Value2 = 0.01322289360846265
$ws.Range("E2").Value2 = 0.0653417645287564
$ws.Range("D3").Value2 = 0.01046569248669804
$ws.Range("E3").Value2 = 0.03228331050675015
$ws.Range("D4").Value2 = 0.01047866112670634
$ws.Range("E4").Value2 = 0.03673788431474723
$ws.Range("D5").Value2 = 0.01134494352726076
$ws.Range("E5").Value2 = 0.04600559527510084
$ws.Range("D6").Value2 = 0.01109842560710299
$ws.Range("E6").Value2 = 0.03025830258302609
$ws.Range("D7").Value2 = 0.01451270448928813
$ws.Range("E7").Value2 = 0.01059001512859314
$ws.Range("D8").Value2 = 0.01093279104699699
$ws.Range("E8").Value2 = 0.02645051194539239
$ws.Range("D9").Value2 = 0.01100559744704358
$ws.Range("E9").Value2 = 0.04209046555858742
$ws.Range("D10").Value2 = 0.01048548672671071
$ws.Range("E10").Value2 = 0.0699127717745085
$ws.Range("D11").Value2 = 0.01112914080712265
$ws.Range("E11").Value2 = 0.04231830726770935
$ws.Range("D12").Value2 = 0.4437732098840148
$ws.Range("E12").Value2 = 0.01190476190476186
$ws.Range("D13").Value2 = 0.01159191648741883
$ws.Range("E13").Value2 = 0.01668335001668342
$ws.Range("D14").Value2 = 0.01074872736687918
$ws.Range("E14").Value2 = 0.01644688101941028
$ws.Range("D15").Value2 = 0.01002589632641657
$ws.Range("E15").Value2 = 0.005491762356465424
$ws.Range("D16").Value2 = 0.009878918406322507
$ws.Range("E16").Value2 = 0.03413173652694601
$ws.Range("D17").Value2 = 0.009491565606074602
$ws.Range("E17").Value2 = 0.02241265655899793
$ws.Range("D18").Value2 = 0.008414599685385344
$ws.Range("E18").Value2 = -0.005597014925373123
$ws.Range("D19").Value2 = 0.009477800646065793
$ws.Range("E19").Value2 = 0.03102719831000789
$ws.Range("D20").Value2 = 0.01059560640678119
$ws.Range("E20").Value2 = -0.001503113592441485
$ws.Range("D21").Value2 = 0.0113728147272786
$ws.Range("E21").Value2 = 0.0437522506301764
$ws.Range("D22").Value2 = 0.01198757376767205
$ws.Range("E22").Value2 = 0.005067567567567766
$ws.Range("D23").Value2 = 0.01113426000712593
$ws.Range("E23").Value2 = 0.02681992337164751
$ws.Range("D24").Value2 = 0.0120392208077051
$ws.Range("E24").Value2 = 0.0514031938013797
$ws.Range("D25").Value2 = 0.01159805952742276
$ws.Range("E25").Value2 = 0.05634023854362846
$ws.Range("D26").Value2 = 0.01121582592717813
$ws.Range("E26").Value2 = 0.05160662122687465
$ws.Range("D27").Value2 = 0.01174799520751872
$ws.Range("E27").Value2 = 0.0008908685968818109
$ws.Range("D28").Value2 = 0.01485489456950713
$ws.Range("E28").Value2 = -0.005559767500631807
$ws.Range("D29").Value2 = 0.01153526400738257
$ws.Range("E29").Value2 = -0.01479289940828399
$ws.Range("D30").Value2 = 0.007219892164620731
$ws.Range("E30").Value2 = 0.0282671036460469
$ws.Range("D31").Value2 = 0.004984963203190376
$ws.Range("E31").Value2 = 0.01118210862619806
$ws.Range("D32").Value2 = 0.009395665926013225
$ws.Range("E32").Value2 = 0.005387931034482873
$ws.Range("D33").Value2 = 0.01107157824708581
$ws.Range("E33").Value2 = -0.0008939213349224984
$ws.Range("D34").Value2 = 0.0103821926466446
$ws.Range("E34").Value2 = 0.02494959677419373
$ws.Range("D35").Value2 = 0.01041336288666455
$ws.Range("E35").Value2 = 0.03071948261924007
$ws.Range("D36").Value2 = 0.009846155526301541
$ws.Range("E36").Value2 = -0.03438395415472795
$ws.Range("D37").Value2 = 0.01106975808708464
$ws.Range("E37").Value2 = 0.02515723270440251
$ws.Range("D38").Value2 = 0.01141240320730394
$ws.Range("E38").Value2 = 0.02330043859649145
$ws.Range("D39").Value2 = 0.01419178752908274
$ws.Range("E39").Value2 = 0.03244837758112107
$ws.Range("D40").Value2 = 0.01065157632681701
$ws.Range("E40").Value2 = 0.02272727272727271
$ws.Range("D41").Value2 = 0.01281028608819858
$ws.Range("E41").Value2 = 0.03225348110258608
$ws.Range("D42").Value2 = 0.0112481337671988
$ws.Range("E42").Value2 = 0.02344350499615677
$ws.Range("D43").Value2 = 0.01135620576726797
$ws.Range("E43").Value2 = 0.02987197724039814
$ws.Range("D44").Value2 = 0.01057319568676685
$ws.Range("E44").Value2 = 0.0393466963622866
$ws.Range("D45").Value2 = 0.01140159600729702
$ws.Range("E45").Value2 = 0.04065851833374889
$ws.Range("D46").Value2 = 0.01081925856692433
$ws.Range("E46").Value2 = 0.05181586860976162
$ws.Range("D47").Value2 = 0.01009574496646128
$ws.Range("E47").Value2 = 0.02484619025082835
$ws.Range("D48").Value2 = 0.009337420805975948
$ws.Range("E48").Value2 = 0.007309941520467822
$ws.Range("D49").Value2 = 0.009372231365998228
$ws.Range("E49").Value2 = 0.04438860971524283
$ws.Range("D50").Value2 = 0.009701907846209222
$ws.Range("E50").Value2 = 0.02026171380329256
$ws.Range("D51").Value2 = 0.009314327525961169
$ws.Range("E51").Value2 = 0.0290435653480221
$ws.Range("D52").Value2 = 0.01001543040640988
$ws.Range("E52").Value2 = 0.04725124943207626
$ws.Range("D53").Value2 = 0.008749054085599395
$ws.Range("E53").Value2 = 0.04244031830238715
$ws.Range("D54").Value2 = 0.004247229602718227
$ws.Range("E54").Value2 = 0.03334672559260721
$ws.Range("D55").Value2 = 0.004158155522661219
$ws.Range("E55").Value2 = 0.02298095863427463
$ws.Range("D56").Value2 = 0.9999999999999998
$ws.Range("E56").Value2 = 0.020037402012824

# Restore sheet protection
$ws.Protect()
